# trafo_id -> gridnode_id refactor
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("windfarms")

# Rename the header cell J1 from "trafo_id" to "gridnode_id"
$ws.Range("J1").Value = "gridnode_id"

# Update the selected cell in the sheet view
$ws.Range("G6").Select()
